$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "NameA"
$ws.Range("B2").Value = "AgeA"
$ws.Range("C2").Value = "CityA"
$ws.Range("D2").Value = "StateA"
$ws.Range("E2").Value = "Zip CodeA"
$ws.Range("F2").Value = "NameB"
$ws.Range("G2").Value = "AgeB"
$ws.Range("H2").Value = "CityB"
$ws.Range("I2").Value = "StateB"
$ws.Range("J2").Value = "Zip CodeB"

$ws.Range("J2").Select()
